$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate rows 237-239 (columns B..AB), leaving A (index) and D (date) untouched ---
# New row237 = old row239, new row238 = old row237, new row239 = old row238
$row237 = $ws.Range("B237:AB237").Value2
$row238 = $ws.Range("B238:AB238").Value2
$row239 = $ws.Range("B239:AB239").Value2

$ws.Range("B237:AB237").Value2 = $row239
$ws.Range("B238:AB238").Value2 = $row237
$ws.Range("B239:AB239").Value2 = $row238

# --- Swap rows 315-316 (columns B..AB), leaving A (index) and D (date) untouched ---
$row315 = $ws.Range("B315:AB315").Value2
$row316 = $ws.Range("B316:AB316").Value2

$ws.Range("B315:AB315").Value2 = $row316
$ws.Range("B316:AB316").Value2 = $row315
